$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1, "sum") onto the
# new header cell (H1) so the new column matches the look of the other
# header cells (bold, centered, bordered) and reuses the same style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" column values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0

$excel.CutCopyMode = 0
